$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 data, matching the style pattern used by row 6 (A6 uses the
# "d-mmm" date style, B6 uses the centered style).
$ws.Range("A7").Value = 44600
$ws.Range("A7").NumberFormat = "d-mmm"

$ws.Range("B7").Value = 2
$ws.Range("B7").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B7").VerticalAlignment = -4108    # xlCenter

$ws.Range("C7").Value = "Experimenting with displays, hover to display, click to display"

$ws.Range("C7").Select()
